$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in missing "url" values (column D) for teams that did not yet have
# a schedule-website slug recorded.
$ws.Range("D3").Value  = "chicago-fire"
$ws.Range("D4").Value  = "colorado-rapids"
$ws.Range("D5").Value  = "columbus-crew-sc"
$ws.Range("D6").Value  = "dc-united"
$ws.Range("D9").Value  = "la-galaxy"
$ws.Range("D15").Value = "new-york-red-bulls"
$ws.Range("D20").Value = "san-jose-earthquakes"
$ws.Range("D24").Value = "vancouver-whitecaps-fc"
$ws.Range("D18").Value = "portland-timbers"
$ws.Range("D17").Value = "philadelphia-union"

# Fix the team_no (UTC/PTC) numbering - shift every value in column E
# (rows 2-25) up by 10.
for ($r = 2; $r -le 25; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $current = $cell.Value2
    $cell.Value = $current + 10
}

# Update the saved view/selection to reflect where the user left off.
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("D25").Select()
